# "Add files via upload" - refresh the rolling quantities table:
#   - every existing date in column A advances forward in the
#     (repeating, period-29) cycle of dates/rates
#   - the two "special rate" rows (the ones using the alternate
#     C/G/J values) move along with the cycle as well
#   - a brand new trailing row (row 49) is appended, continuing the
#     date sequence and reusing the "normal" rate set

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone row 48's formatting (date number format/style) into the new row 49
# before we start overwriting values.
$ws.Range("A48").Copy($ws.Range("A49"))

# Rows 2-3: date shift + switch to the "alternate" C/G/J rate set
$ws.Range("A2").Value = 45558
$ws.Range("C2").Value = 0.00004012
$ws.Range("G2").Value = 280.99031254
$ws.Range("J2").Value = 1941.48834923

$ws.Range("A3").Value = 45559
$ws.Range("C3").Value = 0.00004012
$ws.Range("G3").Value = 280.99031254
$ws.Range("J3").Value = 1941.48834923

# Rows 4-30: date shift only
$ws.Range("A4").Value = 45560
$ws.Range("A5").Value = 45561
$ws.Range("A6").Value = 45562
$ws.Range("A7").Value = 45563
$ws.Range("A8").Value = 45564
$ws.Range("A9").Value = 45565
$ws.Range("A10").Value = 45566
$ws.Range("A11").Value = 45567
$ws.Range("A12").Value = 45568
$ws.Range("A13").Value = 45569
$ws.Range("A14").Value = 45570
$ws.Range("A15").Value = 45571
$ws.Range("A16").Value = 45572
$ws.Range("A17").Value = 45573
$ws.Range("A18").Value = 45574
$ws.Range("A19").Value = 45575
$ws.Range("A20").Value = 45576
$ws.Range("A21").Value = 45577
$ws.Range("A22").Value = 45578
$ws.Range("A23").Value = 45579
$ws.Range("A24").Value = 45580
$ws.Range("A25").Value = 45581
$ws.Range("A26").Value = 45582
$ws.Range("A27").Value = 45583
$ws.Range("A28").Value = 45584
$ws.Range("A29").Value = 45585
$ws.Range("A30").Value = 45586

# Rows 31-32: date shift + switch back to the "normal" C/G/J rate set
$ws.Range("A31").Value = 45587
$ws.Range("C31").Value = 0.00170247
$ws.Range("G31").Value = 465.80531254
$ws.Range("J31").Value = 485.38834923

$ws.Range("A32").Value = 45588
$ws.Range("C32").Value = 0.00170247
$ws.Range("G32").Value = 465.80531254
$ws.Range("J32").Value = 485.38834923

# Rows 33-48: date shift only
$ws.Range("A33").Value = 45589
$ws.Range("A34").Value = 45590
$ws.Range("A35").Value = 45591
$ws.Range("A36").Value = 45592
$ws.Range("A37").Value = 45593
$ws.Range("A38").Value = 45594
$ws.Range("A39").Value = 45595
$ws.Range("A40").Value = 45596
$ws.Range("A41").Value = 45597
$ws.Range("A42").Value = 45598
$ws.Range("A43").Value = 45599
$ws.Range("A44").Value = 45600
$ws.Range("A45").Value = 45601
$ws.Range("A46").Value = 45602
$ws.Range("A47").Value = 45603
$ws.Range("A48").Value = 45604

# Brand new row 49, continuing the date sequence with the "normal" rate set
$ws.Range("A49").Value = 45605
$ws.Range("B49").Value = 116.4121952
$ws.Range("C49").Value = 0.00170247
$ws.Range("D49").Value = 0.008850780000000001
$ws.Range("E49").Value = 0.06933635
$ws.Range("F49").Value = 12792.90181321
$ws.Range("G49").Value = 465.80531254
$ws.Range("H49").Value = 0.24
$ws.Range("I49").Value = 1.7904431
$ws.Range("J49").Value = 485.38834923
